$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{
        E = 3; G = 1.297418666666666; H = 3.892256
        K = 3; M = 0.3922993333333333; N = 1.176898
        O = 0.1197958591217032; P = 0.1197958591217032
        Q = 0.5089764779875555; R = 4.580788301888
        S = 0.1197958591217032; T = 0.1197958591217032
    }
    3 = @{
        E = 3; G = 1.297418666666666; H = 3.892256
        K = 3; M = 0.544463; N = 1.633389
        O = 0.1662618498246574; P = 0.1662618498246574
        Q = 0.7063964595093333; R = 6.357568135584
        S = 0.1662618498246574; T = 0.1662618498246574
    }
    4 = @{
        E = 3; G = 1.297418666666666; H = 3.892256
        K = 3; M = 1.882823333333333; N = 5.64847
        O = 0.5749549377882933; P = 0.5749549377882933
        Q = 2.442810138702222; R = 21.98529124832
        S = 0.5749549377882933; T = 0.5749549377882933
    }
    5 = @{
        E = 3; G = 1.297418666666666; H = 3.892256
        K = 3; M = 0.4551463333333333; N = 1.365439
        O = 0.1389873532653461; P = 0.1389873532653461
        Q = 0.5905153489315554; R = 5.314638140383999
        S = 0.1389873532653461; T = 0.1389873532653461
    }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
